$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1702128.6
$ws.Cells.Item(70, 10).Value = 1738.2
$ws.Cells.Item(70, 12).Value = 5214.6
$ws.Cells.Item(70, 14).Value = -5754.6
$ws.Cells.Item(73, 8).Value = 1702128.6
$ws.Cells.Item(73, 10).Value = 1738.2
$ws.Cells.Item(73, 12).Value = 5214.6
$ws.Cells.Item(73, 14).Value = -7086.6
$ws.Cells.Item(74, 8).Value = 4333
$ws.Cells.Item(74, 9).Value = 3999
$ws.Cells.Item(74, 10).Value = 4399.8
$ws.Cells.Item(74, 11).Value = 3999
$ws.Cells.Item(74, 12).Value = 4399.8
$ws.Cells.Item(74, 13).Value = -3063
$ws.Cells.Item(74, 14).Value = -6271.8
$ws.Cells.Item(77, 8).Value = 4333
$ws.Cells.Item(77, 9).Value = 3999
$ws.Cells.Item(77, 10).Value = 4399.8
$ws.Cells.Item(77, 11).Value = 19995
$ws.Cells.Item(77, 12).Value = 21999
$ws.Cells.Item(77, 13).Value = -15315
$ws.Cells.Item(77, 14).Value = -31359
$ws.Cells.Item(86, 8).Value = 8706399
$ws.Cells.Item(86, 9).Value = 7292.2
$ws.Cells.Item(86, 11).Value = 7292.2
$ws.Cells.Item(86, 13).Value = -6169.2
$ws.Cells.Item(89, 8).Value = 8706399
$ws.Cells.Item(89, 9).Value = 7292.2
$ws.Cells.Item(89, 11).Value = 36461
$ws.Cells.Item(89, 13).Value = -30845
$ws.Cells.Item(106, 8).Value = 59361.89
$ws.Cells.Item(106, 9).Value = 3251
$ws.Cells.Item(106, 11).Value = 3251
$ws.Cells.Item(106, 13).Value = -2620
$ws.Cells.Item(113, 8).Value = 58827428
$ws.Cells.Item(113, 9).Value = 100003020
$ws.Cells.Item(113, 10).Value = 5143
$ws.Cells.Item(113, 11).Value = 100003020
$ws.Cells.Item(113, 12).Value = 5143
$ws.Cells.Item(113, 13).Value = -99999766
$ws.Cells.Item(113, 14).Value = -11651
$ws.Cells.Item(125, 8).Value = 2106.3333
$ws.Cells.Item(125, 9).Value = 2659.5
$ws.Cells.Item(125, 11).Value = 23935.5
$ws.Cells.Item(125, 13).Value = -21475.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1965.0938
$ws.Cells.Item(32, 9).Value = 1965.0938
$ws.Cells.Item(32, 11).Value = 1965.0938
$ws.Cells.Item(32, 13).Value = -1678.0938
$ws.Cells.Item(63, 8).Value = 400000830
$ws.Cells.Item(63, 9).Value = 500001250
$ws.Cells.Item(63, 10).Value = 200000000
$ws.Cells.Item(63, 11).Value = 500001250
$ws.Cells.Item(63, 12).Value = 200000000
$ws.Cells.Item(63, 13).Value = -500000564
$ws.Cells.Item(63, 14).Value = -200001372
$ws.Cells.Item(66, 8).Value = 400000830
$ws.Cells.Item(66, 9).Value = 500001250
$ws.Cells.Item(66, 10).Value = 200000000
$ws.Cells.Item(66, 11).Value = 2500006250
$ws.Cells.Item(66, 12).Value = 1000000000
$ws.Cells.Item(66, 13).Value = -2500002818
$ws.Cells.Item(66, 14).Value = -1000006864
$ws.Cells.Item(122, 8).Value = 13892243
$ws.Cells.Item(122, 9).Value = 17546780
$ws.Cells.Item(122, 11).Value = 52640340
$ws.Cells.Item(122, 13).Value = -52637890
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2740.7917
$ws.Cells.Item(94, 9).Value = 2873.524
$ws.Cells.Item(94, 10).Value = 1811.6666
$ws.Cells.Item(94, 11).Value = 2873.524
$ws.Cells.Item(94, 12).Value = 1811.6666
$ws.Cells.Item(94, 13).Value = -2422.524
$ws.Cells.Item(94, 14).Value = -2713.6666
$ws.Cells.Item(99, 8).Value = 798.8929000000001
$ws.Cells.Item(99, 9).Value = 765.5909
$ws.Cells.Item(99, 11).Value = 765.5909
$ws.Cells.Item(99, 13).Value = 732.4091
$ws.Cells.Item(107, 8).Value = 20018690
$ws.Cells.Item(107, 9).Value = 12489.238
$ws.Cells.Item(107, 11).Value = 12489.238
$ws.Cells.Item(107, 13).Value = -10569.238
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 111118550
$ws.Cells.Item(62, 9).Value = 8596.799999999999
$ws.Cells.Item(62, 11).Value = 8596.799999999999
$ws.Cells.Item(62, 13).Value = -7972.799999999999
$ws.Cells.Item(65, 8).Value = 111118550
$ws.Cells.Item(65, 9).Value = 8596.799999999999
$ws.Cells.Item(65, 11).Value = 42984
$ws.Cells.Item(65, 13).Value = -39864
$ws.Cells.Item(107, 8).Value = 1471.4722
$ws.Cells.Item(107, 9).Value = 1375.5
$ws.Cells.Item(107, 10).Value = 1951.3334
$ws.Cells.Item(107, 11).Value = 1375.5
$ws.Cells.Item(107, 12).Value = 1951.3334
$ws.Cells.Item(107, 13).Value = 544.5
$ws.Cells.Item(107, 14).Value = -5791.3334
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2206.0588
$ws.Cells.Item(34, 10).Value = 2628
$ws.Cells.Item(34, 12).Value = 7884
$ws.Cells.Item(34, 14).Value = -8052
$ws.Cells.Item(39, 8).Value = 2645.5417
$ws.Cells.Item(39, 10).Value = 3815.25
$ws.Cells.Item(39, 12).Value = 11445.75
$ws.Cells.Item(39, 14).Value = -12033.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10420.4
$ws.Cells.Item(70, 9).Value = 9319.286
$ws.Cells.Item(70, 10).Value = 12989.667
$ws.Cells.Item(70, 11).Value = 9319.286
$ws.Cells.Item(70, 12).Value = 12989.667
$ws.Cells.Item(70, 13).Value = -9049.286
$ws.Cells.Item(70, 14).Value = -13529.667
$ws.Cells.Item(73, 8).Value = 10420.4
$ws.Cells.Item(73, 9).Value = 9319.286
$ws.Cells.Item(73, 10).Value = 12989.667
$ws.Cells.Item(73, 11).Value = 9319.286
$ws.Cells.Item(73, 12).Value = 12989.667
$ws.Cells.Item(73, 13).Value = -8383.286
$ws.Cells.Item(73, 14).Value = -14861.667
$ws.Cells.Item(80, 8).Value = 2440.8333
$ws.Cells.Item(80, 9).Value = 3187.4
$ws.Cells.Item(80, 11).Value = 3187.4
$ws.Cells.Item(80, 13).Value = -2189.4
$ws.Cells.Item(83, 8).Value = 2440.8333
$ws.Cells.Item(83, 9).Value = 3187.4
$ws.Cells.Item(83, 11).Value = 15937
$ws.Cells.Item(83, 13).Value = -10945
$ws.Cells.Item(97, 8).Value = 629.4706
$ws.Cells.Item(97, 9).Value = 576.53845
$ws.Cells.Item(97, 11).Value = 576.53845
$ws.Cells.Item(97, 13).Value = -80.53845000000001
$ws.Cells.Item(102, 8).Value = 1554.8948
$ws.Cells.Item(102, 9).Value = 1296.4375
$ws.Cells.Item(102, 10).Value = 2933.3333
$ws.Cells.Item(102, 11).Value = 1296.4375
$ws.Cells.Item(102, 12).Value = 2933.3333
$ws.Cells.Item(102, 13).Value = 325.5625
$ws.Cells.Item(102, 14).Value = -6177.3333
$ws.Cells.Item(107, 8).Value = 1683
$ws.Cells.Item(107, 10).Value = 3750
$ws.Cells.Item(107, 12).Value = 3750
$ws.Cells.Item(107, 14).Value = -7590
$ws.Cells.Item(122, 8).Value = 16131206
$ws.Cells.Item(122, 9).Value = 2203.1667
$ws.Cells.Item(122, 11).Value = 6609.500100000001
$ws.Cells.Item(122, 13).Value = -4159.500100000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 808.3077
$ws.Cells.Item(16, 9).Value = 808.3077
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 808.3077
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(16, 13).Value = -638.3077
$ws.Cells.Item(22, 8).Value = 743.125
$ws.Cells.Item(22, 10).Value = 807.5
$ws.Cells.Item(22, 12).Value = 807.5
$ws.Cells.Item(22, 14).Value = -1397.5
$ws.Cells.Item(27, 8).Value = 743.125
$ws.Cells.Item(27, 10).Value = 807.5
$ws.Cells.Item(27, 12).Value = 807.5
$ws.Cells.Item(27, 14).Value = -1021.5
$ws.Cells.Item(68, 8).Value = 18650.5
$ws.Cells.Item(68, 9).Value = 7634
$ws.Cells.Item(68, 11).Value = 7634
$ws.Cells.Item(68, 13).Value = -6885
$ws.Cells.Item(71, 8).Value = 18650.5
$ws.Cells.Item(71, 9).Value = 7634
$ws.Cells.Item(71, 11).Value = 38170
$ws.Cells.Item(71, 13).Value = -34426
$ws.Cells.Item(82, 8).Value = 1229
$ws.Cells.Item(82, 9).Value = 1194.9166
$ws.Cells.Item(82, 11).Value = 1194.9166
$ws.Cells.Item(82, 13).Value = -833.9166
$ws.Cells.Item(85, 8).Value = 1229
$ws.Cells.Item(85, 9).Value = 1194.9166
$ws.Cells.Item(85, 11).Value = 1194.9166
$ws.Cells.Item(85, 13).Value = 53.08339999999998
$ws.Cells.Item(93, 8).Value = 1283.875
$ws.Cells.Item(93, 9).Value = 1254.8
$ws.Cells.Item(93, 10).Value = 1332.3334
$ws.Cells.Item(93, 11).Value = 1254.8
$ws.Cells.Item(93, 12).Value = 1332.3334
$ws.Cells.Item(93, 13).Value = -6.799999999999955
$ws.Cells.Item(93, 14).Value = -3828.3334
$ws.Cells.Item(122, 8).Value = 3167.75
$ws.Cells.Item(122, 10).Value = 4909.4443
$ws.Cells.Item(122, 12).Value = 14728.3329
$ws.Cells.Item(122, 14).Value = -19628.3329
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 7696471.5
$ws.Cells.Item(81, 9).Value = 1267.5883
$ws.Cells.Item(81, 11).Value = 2535.1766
$ws.Cells.Item(81, 13).Value = -1474.1766
$ws.Cells.Item(84, 8).Value = 7696471.5
$ws.Cells.Item(84, 9).Value = 1267.5883
$ws.Cells.Item(84, 11).Value = 12675.883
$ws.Cells.Item(84, 13).Value = -7371.882999999998
$ws.Cells.Item(100, 8).Value = 3603
$ws.Cells.Item(100, 9).Value = 4645.9165
$ws.Cells.Item(100, 11).Value = 9291.833000000001
$ws.Cells.Item(100, 13).Value = -8750.833000000001
$ws.Cells.Item(136, 8).Value = 6065.6665
$ws.Cells.Item(136, 9).Value = 2823
$ws.Cells.Item(136, 11).Value = 8469
$ws.Cells.Item(136, 13).Value = -5919

Write-Output "Applied 204 cell updates"